$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (46075 -> 46076) for every data row (rows 2 through 161).
$ws.Range("C2:C161").Value = 46076
